$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.712.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.008.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.38"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +8.48%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.999.24"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.71%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +13.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.96"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.122"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.505.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.004.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.695.42"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "438.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.11"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.72"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +11.83%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.35"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.68%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0791"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +15.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.67%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.96%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.24"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.58"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "404.22"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.91%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.767.58"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.10%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.46"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.13"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +20.86%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.68"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.42%  "
